# Resize the cost-comparison table on slide 12 ("Table 4") — the table was
# enlarged (dragged bigger) in the editor, so every column got wider and
# every row got taller, and the table's on-slide position shifted as a
# result. Reproduce by resizing each column / row individually (this keeps
# each <a:gridCol>'s a16:colId extension intact) and then nudging the
# shape's own position to its final on-slide location.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$sh = $s.Shapes.Item(5)
$tbl = $sh.Table

$EMU_PER_POINT = 12700.0

# New column widths (EMU), in left-to-right order.
$colWidths = @(758566, 1067303, 1210054, 1761134, 1761134)
for ($i = 1; $i -le $tbl.Columns.Count; $i++) {
    $tbl.Columns.Item($i).Width = $colWidths[$i - 1] / $EMU_PER_POINT
}

# New row heights (EMU), in top-to-bottom order.
$rowHeights = @(774331, 581360, 774331, 774331, 1016103, 1016103)
for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
    $tbl.Rows.Item($i).Height = $rowHeights[$i - 1] / $EMU_PER_POINT
}

# Final on-slide position/size of the table's graphic frame (EMU).
$sh.Left = 978414 / $EMU_PER_POINT
$sh.Top = 1128683 / $EMU_PER_POINT
$sh.Width = 6558191 / $EMU_PER_POINT
$sh.Height = 4936559 / $EMU_PER_POINT
